$wb = $excel.ActiveWorkbook

# --- "loads" sheet: insert two new leading columns (v_nom_kv, s_base_mva)
# and append two trailing columns (g_shunt_pu, b_shunt_pu) ---
$loads = $wb.Worksheets.Item("loads")

# Shift the existing v_nom_pu / p_nom_mw / q_nom_mvar / bus_idx columns
# two places to the right (Copy preserves the original cell type, e.g. the
# text "1.0" stays text instead of being re-interpreted as a number).
$loads.Range("E1:E2").Copy($loads.Range("G1"))
$loads.Range("D1:D2").Copy($loads.Range("F1"))
$loads.Range("C1:C2").Copy($loads.Range("E1"))
$loads.Range("B1:B2").Copy($loads.Range("D1"))

$loads.Range("B1").Value = "v_nom_kv"
$loads.Range("B2").Value = 22

$loads.Range("C1").Value = "s_base_mva"
$loads.Range("C2").Value = 100

$loads.Range("H1").Value = "g_shunt_pu"
$loads.Range("H2").Value = 0

$loads.Range("I1").Value = "b_shunt_pu"
$loads.Range("I2").Value = 0

# --- Active sheet / selection moved from "lines" to "trafos" ---
$lines = $wb.Worksheets.Item("lines")
$lines.Range("A2:I2").Select()

# --- "loads" sheet selection ends on J2 ---
$loads.Range("J2").Select()

# --- trafos becomes the active/selected tab last ---
$trafos = $wb.Worksheets.Item("trafos")
$trafos.Activate()
$trafos.Range("F16").Select()
